$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1925925925925926
$ws.Range("C2").Value = 0.5703703703703704
$ws.Range("J2").Value = 0.01481481481481482
$ws.Range("P2").Value = 0.1259259259259259
$ws.Range("S2").Value = 0.0962962962962963

# Row 3
$ws.Range("B3").Value = 0.01298701298701299
$ws.Range("C3").Value = 0.01298701298701299
$ws.Range("J3").Value = 0.02597402597402598
$ws.Range("P3").Value = 0.7597402597402597
$ws.Range("S3").Value = 0.1883116883116883

# Row 4
$ws.Range("J4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2222222222222222

# Row 6
$ws.Range("B6").Value = 0.0593607305936073
$ws.Range("D6").Value = 0.0091324200913242
$ws.Range("F6").Value = 0.0365296803652968
$ws.Range("J6").Value = 0.2648401826484018
$ws.Range("O6").Value = 0.0136986301369863
$ws.Range("Q6").Value = 0.1552511415525114
$ws.Range("R6").Value = 0.1004566210045662
$ws.Range("S6").Value = 0.3607305936073059

# Row 7
$ws.Range("B7").Value = 0.072992700729927
$ws.Range("D7").Value = 0.0218978102189781
$ws.Range("E7").Value = 0.0072992700729927
$ws.Range("F7").Value = 0.08759124087591241
$ws.Range("J7").Value = 0.1605839416058394
$ws.Range("O7").Value = 0.0291970802919708
$ws.Range("Q7").Value = 0.1897810218978102
$ws.Range("R7").Value = 0.0948905109489051
$ws.Range("S7").Value = 0.3357664233576642

# Row 8
$ws.Range("B8").Value = 0.08851674641148326
$ws.Range("D8").Value = 0.02392344497607655
$ws.Range("F8").Value = 0.06220095693779904
$ws.Range("J8").Value = 0.1220095693779904
$ws.Range("O8").Value = 0.01913875598086124
$ws.Range("Q8").Value = 0.1626794258373206
$ws.Range("R8").Value = 0.1411483253588517
$ws.Range("S8").Value = 0.3803827751196172

# Row 9
$ws.Range("B9").Value = 0.08823529411764706
$ws.Range("D9").Value = 0.02941176470588235
$ws.Range("F9").Value = 0.07843137254901961
$ws.Range("J9").Value = 0.1176470588235294
$ws.Range("O9").Value = 0.009803921568627451
$ws.Range("Q9").Value = 0.2009803921568628
$ws.Range("R9").Value = 0.09313725490196079
$ws.Range("S9").Value = 0.3823529411764706

# Row 10
$ws.Range("B10").Value = 0.1039370078740157
$ws.Range("D10").Value = 0.01181102362204724
$ws.Range("E10").Value = 0.001574803149606299
$ws.Range("F10").Value = 0.07716535433070866
$ws.Range("J10").Value = 0.1094488188976378
$ws.Range("O10").Value = 0.01259842519685039
$ws.Range("Q10").Value = 0.2188976377952756
$ws.Range("R10").Value = 0.1007874015748031
$ws.Range("S10").Value = 0.3637795275590551

# Row 11
$ws.Range("G11").Value = 0.1441048034934498
$ws.Range("J11").Value = 0.1222707423580786
$ws.Range("K11").Value = 0.2358078602620087
$ws.Range("L11").Value = 0.4759825327510917
$ws.Range("S11").Value = 0.02183406113537118

# Row 12
$ws.Range("G12").Value = 0.7064220183486238
$ws.Range("J12").Value = 0.2293577981651376
$ws.Range("L12").Value = 0.01834862385321101
$ws.Range("S12").Value = 0.04587155963302753

# Row 13
$ws.Range("G13").Value = 0.6071428571428571
$ws.Range("J13").Value = 0.3214285714285715
$ws.Range("S13").Value = 0.07142857142857142

# Row 15
$ws.Range("F15").Value = 0.01463414634146342
$ws.Range("H15").Value = 0.1658536585365854
$ws.Range("I15").Value = 0.1170731707317073
$ws.Range("J15").Value = 0.4146341463414634
$ws.Range("K15").Value = 0.06829268292682927
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.02926829268292683
$ws.Range("S15").Value = 0.1658536585365854

# Row 16
$ws.Range("F16").Value = 0.01136363636363636
$ws.Range("H16").Value = 0.2159090909090909
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.375
$ws.Range("K16").Value = 0.07954545454545454
$ws.Range("M16").Value = 0.02840909090909091
$ws.Range("N16").Value = 0.005681818181818182
$ws.Range("O16").Value = 0.05113636363636364
$ws.Range("S16").Value = 0.1704545454545454

# Row 17
$ws.Range("F17").Value = 0.01805869074492099
$ws.Range("H17").Value = 0.1693002257336343
$ws.Range("I17").Value = 0.0835214446952596
$ws.Range("J17").Value = 0.4785553047404063
$ws.Range("K17").Value = 0.07223476297968397
$ws.Range("M17").Value = 0.01354401805869074
$ws.Range("O17").Value = 0.06094808126410835
$ws.Range("S17").Value = 0.1038374717832957

# Row 18
$ws.Range("F18").Value = 0.01652892561983471
$ws.Range("H18").Value = 0.1859504132231405
$ws.Range("I18").Value = 0.1033057851239669
$ws.Range("J18").Value = 0.4173553719008264
$ws.Range("K18").Value = 0.07851239669421488
$ws.Range("M18").Value = 0.02892561983471074
$ws.Range("O18").Value = 0.05785123966942149
$ws.Range("S18").Value = 0.1115702479338843

# Row 19
$ws.Range("F19").Value = 0.01779359430604982
$ws.Range("H19").Value = 0.201067615658363
$ws.Range("I19").Value = 0.09608540925266904
$ws.Range("J19").Value = 0.400355871886121
$ws.Range("K19").Value = 0.08451957295373666
$ws.Range("M19").Value = 0.02846975088967971
$ws.Range("N19").Value = 0.0008896797153024911
$ws.Range("O19").Value = 0.07829181494661921
$ws.Range("S19").Value = 0.09252669039145907

Write-Host "Updated transition probability matrix with additional simulated games"